$wb = $excel.ActiveWorkbook

# --- Insert new sheet "Feed in from Type 4" after "Feed in from Type 3" ---
$wsType3 = $wb.Worksheets.Item("Feed in from Type 3")
$wsType3.Range("A1:Y4").Copy()
$wsType4 = $wb.Worksheets.Add($null, $wsType3)
$wsType4.Name = "Feed in from Type 4"
$wsType4.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row 1: 0..23
$wsType4.Range("B1").Value = 1
$wsType4.Range("C1").Value = 2
$wsType4.Range("D1").Value = 3
$wsType4.Range("E1").Value = 4
$wsType4.Range("F1").Value = 5
$wsType4.Range("G1").Value = 6
$wsType4.Range("H1").Value = 7
$wsType4.Range("I1").Value = 8
$wsType4.Range("J1").Value = 9
$wsType4.Range("K1").Value = 10
$wsType4.Range("L1").Value = 11
$wsType4.Range("M1").Value = 12
$wsType4.Range("N1").Value = 13
$wsType4.Range("O1").Value = 14
$wsType4.Range("P1").Value = 15
$wsType4.Range("Q1").Value = 16
$wsType4.Range("R1").Value = 17
$wsType4.Range("S1").Value = 18
$wsType4.Range("T1").Value = 19
$wsType4.Range("U1").Value = 20
$wsType4.Range("V1").Value = 21
$wsType4.Range("W1").Value = 22
$wsType4.Range("X1").Value = 23
$wsType4.Range("Y1").Value = 24

# Row 2 data (index 0)
$wsType4.Range("A2").Value = 0
$wsType4.Range("B2").Value = 0
$wsType4.Range("C2").Value = 0
$wsType4.Range("D2").Value = 0
$wsType4.Range("E2").Value = 0
$wsType4.Range("F2").Value = 0
$wsType4.Range("G2").Value = 8
$wsType4.Range("H2").Value = 23
$wsType4.Range("I2").Value = 38
$wsType4.Range("J2").Value = 47
$wsType4.Range("K2").Value = 48
$wsType4.Range("L2").Value = 43
$wsType4.Range("M2").Value = 32
$wsType4.Range("N2").Value = 19
$wsType4.Range("O2").Value = 7
$wsType4.Range("P2").Value = 0
$wsType4.Range("Q2").Value = 0
$wsType4.Range("R2").Value = 0
$wsType4.Range("S2").Value = 0
$wsType4.Range("T2").Value = 0
$wsType4.Range("U2").Value = 0
$wsType4.Range("V2").Value = 0
$wsType4.Range("W2").Value = 0
$wsType4.Range("X2").Value = 0
$wsType4.Range("Y2").Value = 0

# Row 3 data (index 1)
$wsType4.Range("A3").Value = 1
$wsType4.Range("B3").Value = 0
$wsType4.Range("C3").Value = 0
$wsType4.Range("D3").Value = 0
$wsType4.Range("E3").Value = 0
$wsType4.Range("F3").Value = 0
$wsType4.Range("G3").Value = 0
$wsType4.Range("H3").Value = 0
$wsType4.Range("I3").Value = 0
$wsType4.Range("J3").Value = 0
$wsType4.Range("K3").Value = 0
$wsType4.Range("L3").Value = 0
$wsType4.Range("M3").Value = 0
$wsType4.Range("N3").Value = 0
$wsType4.Range("O3").Value = 0
$wsType4.Range("P3").Value = 0
$wsType4.Range("Q3").Value = 0
$wsType4.Range("R3").Value = 0
$wsType4.Range("S3").Value = 0
$wsType4.Range("T3").Value = 0
$wsType4.Range("U3").Value = 0
$wsType4.Range("V3").Value = 0
$wsType4.Range("W3").Value = 0
$wsType4.Range("X3").Value = 0
$wsType4.Range("Y3").Value = 0

# Row 4 data (index 2)
$wsType4.Range("A4").Value = 2
$wsType4.Range("B4").Value = 0
$wsType4.Range("C4").Value = 0
$wsType4.Range("D4").Value = 0
$wsType4.Range("E4").Value = 0
$wsType4.Range("F4").Value = 5
$wsType4.Range("G4").Value = 0
$wsType4.Range("H4").Value = 0
$wsType4.Range("I4").Value = 0
$wsType4.Range("J4").Value = 0
$wsType4.Range("K4").Value = 0
$wsType4.Range("L4").Value = 0
$wsType4.Range("M4").Value = 0
$wsType4.Range("N4").Value = 0
$wsType4.Range("O4").Value = 0
$wsType4.Range("P4").Value = 0
$wsType4.Range("Q4").Value = 0
$wsType4.Range("R4").Value = 0
$wsType4.Range("S4").Value = 0
$wsType4.Range("T4").Value = 0
$wsType4.Range("U4").Value = 0
$wsType4.Range("V4").Value = 0
$wsType4.Range("W4").Value = 0
$wsType4.Range("X4").Value = 0
$wsType4.Range("Y4").Value = 0

# --- DG Dispatch ---
$wsDG = $wb.Worksheets.Item("DG Dispatch")
$wsDG.Range("B2").Value = 120.0100000000004
$wsDG.Range("C2").Value = 120.0099999999948
$wsDG.Range("D2").Value = 120.0099999999948
$wsDG.Range("E2").Value = 120.0100000002157
$wsDG.Range("F2").Value = 0
$wsDG.Range("G2").Value = 0
$wsDG.Range("H2").Value = 0
$wsDG.Range("I2").Value = 0
$wsDG.Range("J2").Value = 0
$wsDG.Range("K2").Value = 0
$wsDG.Range("L2").Value = 0
$wsDG.Range("M2").Value = 0
$wsDG.Range("N2").Value = 0
$wsDG.Range("O2").Value = 120.0099999967206
$wsDG.Range("P2").Value = 240
$wsDG.Range("Q2").Value = 0
$wsDG.Range("R2").Value = 0
$wsDG.Range("S2").Value = 120.0099999999948
$wsDG.Range("T2").Value = 0
$wsDG.Range("U2").Value = 160.7883227884676
$wsDG.Range("V2").Value = 120.0099999999948
$wsDG.Range("W2").Value = 0
$wsDG.Range("X2").Value = 240
$wsDG.Range("S4").Value = 120.4749999999849

# --- PV Dispatch ---
$wsPV = $wb.Worksheets.Item("PV Dispatch")
$wsPV.Range("G2").Value = 683.8020000000001
$wsPV.Range("H2").Value = 1918.728000000001
$wsPV.Range("I2").Value = 3163.86
$wsPV.Range("J2").Value = 3939.516000000001
$wsPV.Range("K2").Value = 4051.782000000001
$wsPV.Range("L2").Value = 3572.1
$wsPV.Range("M2").Value = 2684.178000000001
$wsPV.Range("N2").Value = 1612.548
$wsPV.Range("O2").Value = 642.9780000000001
$wsPV.Range("P2").Value = 20.412

# --- Battery Input  ---
$wsBatIn = $wb.Worksheets.Item("Battery Input ")
$wsBatIn.Range("H2").Value = 1165.948
$wsBatIn.Range("I2").Value = 2635.190000003213
$wsBatIn.Range("J2").Value = 3595.936000000801
$wsBatIn.Range("K2").Value = 3631.47700000079
$wsBatIn.Range("L2").Value = 3055.229999998706
$wsBatIn.Range("M2").Value = 2016.553000002045
$wsBatIn.Range("N2").Value = 710.8380000014528
$wsBatIn.Range("T2").Value = 0.00000002596061676740646
$wsBatIn.Range("I3").Value = 0
$wsBatIn.Range("K3").Value = 5305
$wsBatIn.Range("L3").Value = 2947.800306697911
$wsBatIn.Range("N3").Value = 0
$wsBatIn.Range("O3").Value = 3332.994000000002
$wsBatIn.Range("P3").Value = 0
$wsBatIn.Range("V3").Value = -0.00000004435423761606216
$wsBatIn.Range("B4").Value = -0.00000001653097569942474
$wsBatIn.Range("H4").Value = 0
$wsBatIn.Range("K4").Value = 5305
$wsBatIn.Range("L4").Value = 0
$wsBatIn.Range("M4").Value = 5305
$wsBatIn.Range("N4").Value = 3870.642747474748
$wsBatIn.Range("O4").Value = 2666.832
$wsBatIn.Range("V4").Value = 0.000000016065314412117

# --- Battery Output ---
$wsBatOut = $wb.Worksheets.Item("Battery Output")
$wsBatOut.Range("B2").Value = 1091.5
$wsBatOut.Range("C2").Value = 1016.25
$wsBatOut.Range("D2").Value = 941.8
$wsBatOut.Range("E2").Value = 911.1599999997844
$wsBatOut.Range("F2").Value = 977.6299999999137
$wsBatOut.Range("G2").Value = 175.0929999981854
$wsBatOut.Range("O2").Value = 348.4570000002062
$wsBatOut.Range("P2").Value = 900.4780000012869
$wsBatOut.Range("Q2").Value = 1191.31
$wsBatOut.Range("R2").Value = 1141.460000000069
$wsBatOut.Range("T2").Value = 1289.760000025961
$wsBatOut.Range("U2").Value = 1218.411677211533
$wsBatOut.Range("V2").Value = 1226.740000000084
$wsBatOut.Range("W2").Value = 1297.99
$wsBatOut.Range("X2").Value = 900.13
$wsBatOut.Range("Y2").Value = 824.4600000000564
$wsBatOut.Range("E3").Value = 345.4699999959026
$wsBatOut.Range("Q3").Value = 245.1309999986263
$wsBatOut.Range("R3").Value = 891.735999999749
$wsBatOut.Range("V3").Value = 1172.999999958121
$wsBatOut.Range("B4").Value = 1150.349999987096
$wsBatOut.Range("C4").Value = 1086.22
$wsBatOut.Range("E4").Value = 1001.979999999927
$wsBatOut.Range("F4").Value = 396.9649999997055
$wsBatOut.Range("Q4").Value = 1298.34
$wsBatOut.Range("R4").Value = 1258.219999999903
$wsBatOut.Range("S4").Value = 1148.635000000015
$wsBatOut.Range("U4").Value = 1520.979999999277
$wsBatOut.Range("V4").Value = 1498.260000016065

# --- State of Charge ---
$wsSoC = $wb.Worksheets.Item("State of Charge")
$wsSoC.Range("B2").Value = 8306.558585857087
$wsSoC.Range("C2").Value = 7280.043434341937
$wsSoC.Range("D2").Value = 6328.7303030288
$wsSoC.Range("E2").Value = 5408.366666665382
$wsSoC.Range("F2").Value = 4420.861616160101
$wsSoC.Range("G2").Value = 4244
$wsSoC.Range("H2").Value = 5398.28852
$wsSoC.Range("I2").Value = 8007.126620003181
$wsSoC.Range("J2").Value = 11567.10326000397
$wsSoC.Range("K2").Value = 15162.26549000476
$wsSoC.Range("L2").Value = 18186.94319000348
$wsSoC.Range("M2").Value = 20183.33066000843
$wsSoC.Range("N2").Value = 20887.06028000987
$wsSoC.Range("O2").Value = 20535.08351233289
$wsSoC.Range("P2").Value = 19625.50977495827
$wsSoC.Range("Q2").Value = 18422.16634061527
$wsSoC.Range("R2").Value = 17269.17644162582
$wsSoC.Range("S2").Value = 16234.83300728341
$wsSoC.Range("T2").Value = 14932.04512849501
$wsSoC.Range("U2").Value = 13701.32626262477
$wsSoC.Range("V2").Value = 12462.19494949346
$wsSoC.Range("W2").Value = 11151.09393939245
$wsSoC.Range("X2").Value = 10241.87171717022
$wsSoC.Range("Y2").Value = 9409.083838382343
$wsSoC.Range("B3").Value = 6163.161616158211
$wsSoC.Range("C3").Value = 5343.969696966286
$wsSoC.Range("D3").Value = 4592.959595955797
$wsSoC.Range("E3").Value = 4244
$wsSoC.Range("F3").Value = 4244
$wsSoC.Range("G3").Value = 4244
$wsSoC.Range("H3").Value = 4244
$wsSoC.Range("I3").Value = 4244
$wsSoC.Range("J3").Value = 4244
$wsSoC.Range("K3").Value = 9495.95
$wsSoC.Range("L3").Value = 12414.27230363097
$wsSoC.Range("M3").Value = 12414.27230363097
$wsSoC.Range("N3").Value = 12414.27230363097
$wsSoC.Range("O3").Value = 15713.93636363097
$wsSoC.Range("P3").Value = 15713.93636363097
$wsSoC.Range("Q3").Value = 15466.32929292494
$wsSoC.Range("R3").Value = 14565.58585858176
$wsSoC.Range("S3").Value = 13605.18181817671
$wsSoC.Range("T3").Value = 12490.53535353032
$wsSoC.Range("U3").Value = 11279.92929292426
$wsSoC.Range("V3").Value = 10095.08080807667
$wsSoC.Range("W3").Value = 8953.868686864615
$wsSoC.Range("X3").Value = 7966.494949491539
$wsSoC.Range("Y3").Value = 7051.343434340024
$wsSoC.Range("B4").Value = 7786.08585858586
$wsSoC.Range("C4").Value = 6688.89393939394
$wsSoC.Range("D4").Value = 5657.075757575757
$wsSoC.Range("E4").Value = 4644.974747474622
$wsSoC.Range("F4").Value = 4244
$wsSoC.Range("G4").Value = 4244
$wsSoC.Range("H4").Value = 4244
$wsSoC.Range("I4").Value = 4244
$wsSoC.Range("J4").Value = 4244
$wsSoC.Range("K4").Value = 9495.95
$wsSoC.Range("L4").Value = 9495.95
$wsSoC.Range("M4").Value = 14747.9
$wsSoC.Range("N4").Value = 18579.83632
$wsSoC.Range("O4").Value = 21220
$wsSoC.Range("P4").Value = 21220
$wsSoC.Range("Q4").Value = 19908.54545454546
$wsSoC.Range("R4").Value = 18637.61616161616
$wsSoC.Range("S4").Value = 17477.37878787878
$wsSoC.Range("T4").Value = 16038.39898989897
$wsSoC.Range("U4").Value = 14502.05555555555
$wsSoC.Range("V4").Value = 12988.66161616162
$wsSoC.Range("W4").Value = 11517.5
$wsSoC.Range("X4").Value = 10197.36868686835
$wsSoC.Range("Y4").Value = 8948.055555555555

# --- Feed in from Type 2 ---
$wsType2 = $wb.Worksheets.Item("Feed in from Type 2")
$wsType2.Range("M2").Value = 116.5350000000001

# --- Feed in from Type 3 ---
$wsType3b = $wb.Worksheets.Item("Feed in from Type 3")
$wsType3b.Range("H2").Value = 49.03999999999998

# --- Costs and Revenues ---
$wsCosts = $wb.Worksheets.Item("Costs and Revenues")
$wsCosts.Range("B2").Value = 3384782.32645785
$wsCosts.Range("C2").Value = 1200
$wsCosts.Range("D2").Value = 426180.0239132897
$wsCosts.Range("E2").Value = 187615
$wsCosts.Range("F2").Value = 2495303.876166271

# --- Capacities ---
$wsCap = $wb.Worksheets.Item("Capacities")
$wsCap.Range("C3").Value = 10206
$wsCap.Range("B4").Value = 6
$wsCap.Range("C4").Value = 5305

# --- Connected Households ---
$wsHH = $wb.Worksheets.Item("Connected Households")
$wsHH.Range("G3").Value = 510
$wsHH.Range("H3").Value = 510
$wsHH.Range("J3").Value = 509
$wsHH.Range("G4").Value = 40
$wsHH.Range("H4").Value = 40
$wsHH.Range("J4").Value = 40

# --- Yearly demand ---
$wsYD = $wb.Worksheets.Item("Yearly demand")
$wsYD.Range("E2").Value = -1031.17
$wsYD.Range("G2").Value = -858.895
$wsYD.Range("H2").Value = -752.78
$wsYD.Range("I2").Value = -528.6700000000001
$wsYD.Range("J2").Value = -343.5800000000002
$wsYD.Range("K2").Value = -420.3050000000001
$wsYD.Range("L2").Value = -516.8699999999999
$wsYD.Range("M2").Value = -667.625
$wsYD.Range("N2").Value = -901.71
$wsYD.Range("O2").Value = -1111.445
$wsYD.Range("Q2").Value = -1191.31
$wsYD.Range("E3").Value = -651.6500000000001
$wsYD.Range("F3").Value = -316.1400000000001
$wsYD.Range("G3").Value = 69.27000000000005
$wsYD.Range("H3").Value = 254.87
$wsYD.Range("I3").Value = 454.3449999999998
$wsYD.Range("J3").Value = 618.7799999999999
$wsYD.Range("K3").Value = 570.675
$wsYD.Range("L3").Value = 538.2999999999998
$wsYD.Range("M3").Value = 445.57
$wsYD.Range("N3").Value = 215.0700000000001
$wsYD.Range("O3").Value = -112.885
$wsYD.Range("P3").Value = -434.76
$wsYD.Range("Q3").Value = -847.285
$wsYD.Range("E4").Value = -1001.98
$wsYD.Range("F4").Value = -856.235
$wsYD.Range("G4").Value = -396.7049999999999
$wsYD.Range("H4").Value = -85.70999999999998
$wsYD.Range("I4").Value = 172.9800000000001
$wsYD.Range("J4").Value = 364.04
$wsYD.Range("K4").Value = 320.71
$wsYD.Range("L4").Value = 275.8149999999999
$wsYD.Range("M4").Value = 158.8449999999999
$wsYD.Range("N4").Value = -115.765
$wsYD.Range("O4").Value = -530.1600000000001
$wsYD.Range("P4").Value = -976.62
$wsYD.Range("Q4").Value = -1298.34

# Connected Households: add new row 5 for household type "3"
$wsHH.Range("A5").Value = 3
$wsHH.Range("B5").Value = 1
$wsHH.Range("C5").Value = 1
$wsHH.Range("D5").Value = 1
$wsHH.Range("E5").Value = 1
$wsHH.Range("F5").Value = 1
$wsHH.Range("G5").Value = 1
$wsHH.Range("H5").Value = 1
$wsHH.Range("I5").Value = 1
$wsHH.Range("J5").Value = 1
$wsHH.Range("K5").Value = 1
$wsHH.Range("L5").Value = 1
$wsHH.Range("M5").Value = 1
$wsHH.Range("N5").Value = 1
$wsHH.Range("O5").Value = 1
$wsHH.Range("P5").Value = 1

Write-Output "done"